# Auto-generated edit script: update Leve market-board cached values
# across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets (scheduled market refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1184
$ws.Range("I19").Value = 1200
$ws.Range("J19").Value = 1176
$ws.Range("K19").Value = 1200
$ws.Range("L19").Value = 1176
$ws.Range("M19").Value = -1025
$ws.Range("N19").Value = -1526

# Row 40
$ws.Range("H40").Value = 7757.4
$ws.Range("I40").Value = 5214.2856
$ws.Range("K40").Value = 5214.2856
$ws.Range("M40").Value = -5039.2856

# Row 43
$ws.Range("H43").Value = 1230.5
$ws.Range("I43").Value = 129.66667
$ws.Range("J43").Value = 2331.3333
$ws.Range("K43").Value = 129.66667
$ws.Range("L43").Value = 2331.3333
$ws.Range("M43").Value = -60.66667000000001
$ws.Range("N43").Value = -2469.3333

# Row 55
$ws.Range("H55").Value = 395.33334
$ws.Range("I55").Value = 395.18182
$ws.Range("K55").Value = 395.18182
$ws.Range("M55").Value = -181.18182

# Row 132
$ws.Range("H132").Value = 3376.9363
$ws.Range("I132").Value = 3603.4883
$ws.Range("J132").Value = 941.5
$ws.Range("K132").Value = 10810.4649
$ws.Range("L132").Value = 2824.5
$ws.Range("M132").Value = -8280.464899999999
$ws.Range("N132").Value = -7884.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1685.1406
$ws.Range("I32").Value = 1685.1406
$ws.Range("K32").Value = 1685.1406
$ws.Range("M32").Value = -1398.1406

# Row 45
$ws.Range("H45").Value = 48301.89
$ws.Range("I45").Value = 85143.39999999999
$ws.Range("K45").Value = 85143.39999999999
$ws.Range("M45").Value = -84766.39999999999

# Row 61
$ws.Range("H61").Value = 1660.1177
$ws.Range("I61").Value = 1094.0385
$ws.Range("K61").Value = 1094.0385
$ws.Range("M61").Value = -882.0385000000001

# Row 102
$ws.Range("H102").Value = 66670340
$ws.Range("I102").Value = 125002600
$ws.Range("J102").Value = 4899.7144
$ws.Range("K102").Value = 125002600
$ws.Range("L102").Value = 4899.7144
$ws.Range("M102").Value = -125000978
$ws.Range("N102").Value = -8143.7144

# Row 132
$ws.Range("H132").Value = 1839.6046
$ws.Range("I132").Value = 1222.0312
$ws.Range("J132").Value = 3636.182
$ws.Range("K132").Value = 3666.0936
$ws.Range("L132").Value = 10908.546
$ws.Range("M132").Value = -1136.0936
$ws.Range("N132").Value = -15968.546

# Row 136
$ws.Range("H136").Value = 1660.1177
$ws.Range("I136").Value = 1094.0385
$ws.Range("K136").Value = 3282.1155
$ws.Range("M136").Value = -732.1155000000003

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 31255662
$ws.Range("I20").Value = 38467828
$ws.Range("K20").Value = 38467828
$ws.Range("M20").Value = -38467581

# Row 86
$ws.Range("H86").Value = 4299.769
$ws.Range("I86").Value = 3656.8572
$ws.Range("J86").Value = 5049.8335
$ws.Range("K86").Value = 3656.8572
$ws.Range("L86").Value = 5049.8335
$ws.Range("M86").Value = -2533.8572
$ws.Range("N86").Value = -7295.8335

# Row 89
$ws.Range("H89").Value = 4299.769
$ws.Range("I89").Value = 3656.8572
$ws.Range("J89").Value = 5049.8335
$ws.Range("K89").Value = 18284.286
$ws.Range("L89").Value = 25249.1675
$ws.Range("M89").Value = -12668.286
$ws.Range("N89").Value = -36481.1675

# Row 99
$ws.Range("H99").Value = 2534.3572
$ws.Range("I99").Value = 2123.4167
$ws.Range("K99").Value = 2123.4167
$ws.Range("M99").Value = -625.4167000000002

# Row 107
$ws.Range("H107").Value = 3663922
$ws.Range("I107").Value = 5918105.5
$ws.Range("K107").Value = 5918105.5
$ws.Range("M107").Value = -5916185.5

# Row 134
$ws.Range("H134").Value = 1708.4546
$ws.Range("I134").Value = 1298.9788
$ws.Range("J134").Value = 4114.125
$ws.Range("K134").Value = 3896.936400000001
$ws.Range("L134").Value = 12342.375
$ws.Range("M134").Value = -1361.936400000001
$ws.Range("N134").Value = -17412.375

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 1577.7273
$ws.Range("I94").Value = 986
$ws.Range("J94").Value = 1915.8572
$ws.Range("K94").Value = 986
$ws.Range("L94").Value = 1915.8572
$ws.Range("M94").Value = -535
$ws.Range("N94").Value = -2817.8572

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1291.29
$ws.Range("I4").Value = 1307.4387
$ws.Range("K4").Value = 3922.3161
$ws.Range("M4").Value = -3810.3161

# Row 5
$ws.Range("H5").Value = 1330.8334
$ws.Range("J5").Value = 1372
$ws.Range("L5").Value = 4116
$ws.Range("N5").Value = -4340

# Row 135
$ws.Range("H135").Value = 1330.8334
$ws.Range("J135").Value = 1372
$ws.Range("L135").Value = 12348
$ws.Range("N135").Value = -17418

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 70
$ws.Range("H70").Value = 8329.333000000001
$ws.Range("J70").Value = 9994.75
$ws.Range("L70").Value = 9994.75
$ws.Range("N70").Value = -10534.75

# Row 73
$ws.Range("H73").Value = 8329.333000000001
$ws.Range("J73").Value = 9994.75
$ws.Range("L73").Value = 9994.75
$ws.Range("N73").Value = -11866.75

# Row 97
$ws.Range("H97").Value = 1249.3125
$ws.Range("I97").Value = 725.9
$ws.Range("J97").Value = 2121.6667
$ws.Range("K97").Value = 725.9
$ws.Range("L97").Value = 2121.6667
$ws.Range("M97").Value = -229.9
$ws.Range("N97").Value = -3113.6667

# Row 126
$ws.Range("H126").Value = 7127.4443
$ws.Range("I126").Value = 2367.375
$ws.Range("J126").Value = 10935.5
$ws.Range("K126").Value = 7102.125
$ws.Range("L126").Value = 32806.5
$ws.Range("M126").Value = -4632.125
$ws.Range("N126").Value = -37746.5

# Row 132
$ws.Range("H132").Value = 2027.4
$ws.Range("J132").Value = 3810.7778
$ws.Range("L132").Value = 11432.3334
$ws.Range("N132").Value = -16492.3334

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 9551
$ws.Range("I14").Value = 7734.6665
$ws.Range("K14").Value = 7734.6665
$ws.Range("M14").Value = -7566.6665

# Row 62
$ws.Range("H62").Value = 7349.2856
$ws.Range("I62").Value = 6925
$ws.Range("J62").Value = 7519
$ws.Range("K62").Value = 6925
$ws.Range("L62").Value = 7519
$ws.Range("M62").Value = -6301
$ws.Range("N62").Value = -8767

# Row 65
$ws.Range("H65").Value = 7349.2856
$ws.Range("I65").Value = 6925
$ws.Range("J65").Value = 7519
$ws.Range("K65").Value = 34625
$ws.Range("L65").Value = 37595
$ws.Range("M65").Value = -31505
$ws.Range("N65").Value = -43835

# Row 132
$ws.Range("H132").Value = 3702.6155
$ws.Range("I132").Value = 3596.8635
$ws.Range("K132").Value = 10790.5905
$ws.Range("M132").Value = -8260.5905
